# Secdep Loan, Saving, RD scenarios
# Updates the "Summary" and "Repayment schedule" sheets with the results of
# a re-run of the loan scenario, and moves the active tab / selections to
# match where the author ended up after the run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: a couple of totals change because the repayment numbers
# downstream changed.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("A3").Value = 462.45
$wsSummary.Range("E3").Value = 462.45
$wsSummary.Range("F3").Value = 0

# F2 goes to 0 and also picks up F3's (General) number format instead of
# its previous #,##0.00 one, so copy the formatting across before writing
# the new value.
$wsSummary.Range("F3").Copy() | Out-Null
$wsSummary.Range("F2").PasteSpecial(-4122) | Out-Null
$wsSummary.Range("F2").Value = 0

# ---------------------------------------------------------------------
# Repayment schedule sheet: the whole amortisation table is recomputed.
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$scheduleRows = @(
    @{Row=3;  B=90; C=42095; F=788.49; G=9211.51;            H=100},
    @{Row=4;  B=0;  C=42095; F=888.49; G=8323.02;             H=0},
    @{Row=5;  B=0;            F=888.49; G=7434.53;             H=0},
    @{Row=6;            F=814.14; G=6620.39;             H=74.349999999999994},
    @{Row=7;            F=822.29; G=5798.1;              H=66.2},
    @{Row=8;            F=830.51; G=4967.59;             H=57.98},
    @{Row=9;            F=835.61; G=4131.9799999999996;  H=52.88},
    @{Row=10;           F=841.52; G=3290.46;              H=46.97},
    @{Row=11;           F=855.59; G=2434.87;              H=32.9},
    @{Row=12; B=61; C=42339; F=864.14; G=1570.73;          H=24.35},
    @{Row=13; B=0;            F=888.49; G=682.24;           H=0},
    @{Row=14;           F=682.24;                           H=6.82}
)

foreach ($r in $scheduleRows) {
    if ($r.ContainsKey("B")) { $wsRepay.Cells.Item($r.Row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $wsRepay.Cells.Item($r.Row, 3).Value = $r.C }
    if ($r.ContainsKey("F")) { $wsRepay.Cells.Item($r.Row, 6).Value = $r.F }
    if ($r.ContainsKey("G")) { $wsRepay.Cells.Item($r.Row, 7).Value = $r.G }
    if ($r.ContainsKey("H")) { $wsRepay.Cells.Item($r.Row, 8).Value = $r.H }
}

# Row 14's "Due"/"Over Due" totals (K14 / Q14) also change.
$wsRepay.Cells.Item(14, 11).Value = 689.06
$wsRepay.Cells.Item(14, 17).Value = 689.06

# ---------------------------------------------------------------------
# Selections / active tab: the author ended up with C2:D2 selected on
# Summary and then moved on to (and left the workbook on) the Repayment
# schedule sheet with J15 selected, so Repayment schedule must be the
# last sheet activated/selected.
# ---------------------------------------------------------------------
$wsSummary.Range("C2:D2").Select()

$wsRepay.Range("J15").Select()
